# ItemTable.xlsx edit: split StatusComponent into Box-related sheets, add BoxComponent data, tweak UI selection.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert two new worksheets ("BoxData", "BoxStarting") right before "StringKey",
#    in that final order (BoxData, then BoxStarting, then StringKey).
# ---------------------------------------------------------------------------
$itemStartingPlayer = $wb.Worksheets.Item("ItemStartingPlayer")
$boxData = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $itemStartingPlayer)
$boxData.Name = "BoxData"
$boxStarting = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $boxData)
$boxStarting.Name = "BoxStarting"

# ---------------------------------------------------------------------------
# 2. Populate "BoxStarting" (boxId, areaIndex) first so new shared strings are
#    allocated in the same order as the source workbook (boxId, areaIndex, itemId).
# ---------------------------------------------------------------------------
$boxStarting.Range("A1").Value = "int"
$boxStarting.Range("B1").Value = "int"
$boxStarting.Range("A2").Value = "boxId"
$boxStarting.Range("B2").Value = "areaIndex"

$boxStartingData = New-Object 'object[,]' 2,2
$boxStartingData[0,0] = 1
$boxStartingData[0,1] = 1
$boxStartingData[1,0] = 2
$boxStartingData[1,1] = 2
$boxStarting.Range("A3:B4").Value = $boxStartingData

$boxStarting.Columns.Item(1).ColumnWidth = 12.38
$boxStarting.Columns.Item(2).ColumnWidth = 16.68

$boxStarting.Activate()
$boxStarting.Range("D6").Select()

# ---------------------------------------------------------------------------
# 3. Populate "BoxData" (itemId, boxId).
# ---------------------------------------------------------------------------
$boxData.Range("A1").Value = "int"
$boxData.Range("B1").Value = "int"
$boxData.Range("A2").Value = "itemId"
$boxData.Range("B2").Value = "boxId"

$boxDataRows = New-Object 'object[,]' 12,2
$boxDataRows[0,0] = 201
$boxDataRows[0,1] = 1
$boxDataRows[1,0] = 202
$boxDataRows[1,1] = 1
$boxDataRows[2,0] = 203
$boxDataRows[2,1] = 1
$boxDataRows[3,0] = 204
$boxDataRows[3,1] = 1
$boxDataRows[4,0] = 205
$boxDataRows[4,1] = 1
$boxDataRows[5,0] = 206
$boxDataRows[5,1] = 1
$boxDataRows[6,0] = 301
$boxDataRows[6,1] = 2
$boxDataRows[7,0] = 302
$boxDataRows[7,1] = 2
$boxDataRows[8,0] = 303
$boxDataRows[8,1] = 2
$boxDataRows[9,0] = 304
$boxDataRows[9,1] = 2
$boxDataRows[10,0] = 305
$boxDataRows[10,1] = 2
$boxDataRows[11,0] = 306
$boxDataRows[11,1] = 2
$boxData.Range("A3:B14").Value = $boxDataRows

$boxData.Columns.Item(1).ColumnWidth = 19.48
$boxData.Columns.Item(2).ColumnWidth = 20.98
$boxData.PageSetup.PaperSize = 9
$boxData.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Append the new item-type rows (201-206, 301-306) to "ItemType".
# ---------------------------------------------------------------------------
$itemType = $wb.Worksheets.Item("ItemType")
$itemTypeRows = New-Object 'object[,]' 12,2
$itemTypeRows[0,0] = 201
$itemTypeRows[0,1] = "Armor"
$itemTypeRows[1,0] = 202
$itemTypeRows[1,1] = "Armor"
$itemTypeRows[2,0] = 203
$itemTypeRows[2,1] = "Smoke"
$itemTypeRows[3,0] = 204
$itemTypeRows[3,1] = "Smoke"
$itemTypeRows[4,0] = 205
$itemTypeRows[4,1] = "Rader"
$itemTypeRows[5,0] = 206
$itemTypeRows[5,1] = "Rader"
$itemTypeRows[6,0] = 301
$itemTypeRows[6,1] = "Armor"
$itemTypeRows[7,0] = 302
$itemTypeRows[7,1] = "Armor"
$itemTypeRows[8,0] = 303
$itemTypeRows[8,1] = "Teleport"
$itemTypeRows[9,0] = 304
$itemTypeRows[9,1] = "Teleport"
$itemTypeRows[10,0] = 305
$itemTypeRows[10,1] = "Teleport"
$itemTypeRows[11,0] = 306
$itemTypeRows[11,1] = "Teleport"
$itemType.Range("A17:B28").Value = $itemTypeRows

# ---------------------------------------------------------------------------
# 5. UI bookkeeping: selections + active sheet/tab to match the saved workbook state.
# ---------------------------------------------------------------------------
$itemType.Range("B22").Select()

$itemPosition = $wb.Worksheets.Item("ItemPosition")
# (ItemPosition loses tabSelected automatically once another sheet is activated below.)

$boxData.Activate()

Write-Host "done"
